$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 142 (shifting the existing
# rows 142-149 down to 144-151) to make room for a new week of data.
$ws.Rows("142:143").Insert()

# New row 142 — Coliflor, Segunda, week of 45041
$ws.Range("A142").Value = 1
$ws.Range("B142").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C142").Value = "Arica y Parinacota"
$ws.Range("D142").Value = 45041
$ws.Range("E142").Value = 15
$ws.Range("F142").Value = 100112008
$ws.Range("G142").Value = "Coliflor"
$ws.Range("H142").Value = "Sin especificar"
$ws.Range("I142").Value = "Segunda"
$ws.Range("J142").Value = 800
$ws.Range("K142").Value = 1300
$ws.Range("L142").Value = 1400
$ws.Range("M142").Value = 1350
$ws.Range("N142").Value = "$/unidad"
$ws.Range("O142").Value = "Región de Arica y Parinacota"
$ws.Range("P142").Value = 1350
$ws.Range("Q142").Value = 1
$ws.Range("R142").Value = "Hortaliza"

# New row 143 — Coliflor, Tercera, week of 45041
$ws.Range("A143").Value = 1
$ws.Range("B143").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C143").Value = "Arica y Parinacota"
$ws.Range("D143").Value = 45041
$ws.Range("E143").Value = 15
$ws.Range("F143").Value = 100112008
$ws.Range("G143").Value = "Coliflor"
$ws.Range("H143").Value = "Sin especificar"
$ws.Range("I143").Value = "Tercera"
$ws.Range("J143").Value = 1000
$ws.Range("K143").Value = 900
$ws.Range("L143").Value = 1000
$ws.Range("M143").Value = 950
$ws.Range("N143").Value = "$/unidad"
$ws.Range("O143").Value = "Región de Arica y Parinacota"
$ws.Range("P143").Value = 950
$ws.Range("Q143").Value = 1
$ws.Range("R143").Value = "Hortaliza"
